$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the numeric-looking Price/Volume cells that change,
# so Excel does not auto-convert strings like "309.30" or "2.06%" into numbers.
$numericCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "D6", "E6", "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "D15", "E15", "D16", "E16", "D17", "E17", "E18", "E19", "E20", "D21", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D39", "E39", "D40", "E40", "D41", "E41", "E42", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "E47", "E48", "E49", "E50", "E51")
foreach ($cellRef in $numericCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated values cell by cell, matching the target diff.
$ws.Range("D2").Value = "309.30"
$ws.Range("E2").Value = "2.06%"
$ws.Range("D3").Value = "38.93"
$ws.Range("E3").Value = "9.01%"
$ws.Range("D4").Value = "5.080"
$ws.Range("E4").Value = "0.80%"
$ws.Range("D5").Value = "0.08194"
$ws.Range("D6").Value = "2.019"
$ws.Range("E6").Value = "9.46%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "4.176"
$ws.Range("E7").Value = "1.87%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "7.911"
$ws.Range("E8").Value = "1.62%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "0.9319"
$ws.Range("E9").Value = "1.33%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "0.1416"
$ws.Range("E10").Value = "4.63%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1954"
$ws.Range("E11").Value = "3.22%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.09364"
$ws.Range("E12").Value = "3.61%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "0.03461"
$ws.Range("E13").Value = "-0.28%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "0.09840"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "0.001414"
$ws.Range("E15").Value = "0.94%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "0.005973"
$ws.Range("E16").Value = "-2.33%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "3.711"
$ws.Range("E17").Value = "-0.32%"
$ws.Range("E18").Value = "4.47%"
$ws.Range("E19").Value = "0.41%"
$ws.Range("E20").Value = "-0.53%"
$ws.Range("D21").Value = "4.817"
$ws.Range("E21").Value = "-6.53%"
$ws.Range("D22").Value = "0.2452"
$ws.Range("E22").Value = "11.84%"
$ws.Range("D23").Value = "0.04471"
$ws.Range("E23").Value = "1.56%"
$ws.Range("D24").Value = "0.001240"
$ws.Range("E24").Value = "0.34%"
$ws.Range("E25").Value = "-9.44%"
$ws.Range("D39").Value = "0.02136"
$ws.Range("E39").Value = "10.36%"
$ws.Range("D40").Value = "0.05191"
$ws.Range("E40").Value = "0.62%"
$ws.Range("D41").Value = "0.007481"
$ws.Range("E41").Value = "-1.77%"
$ws.Range("E42").Value = "-0.41%"
$ws.Range("E43").Value = "2.16%"
$ws.Range("D44").Value = "0.002132"
$ws.Range("E44").Value = "-1.35%"
$ws.Range("D45").Value = "0.009683"
$ws.Range("E45").Value = "-4.75%"
$ws.Range("D46").Value = "0.00006298"
$ws.Range("E46").Value = "2.42%"
$ws.Range("E47").Value = "-0.01%"
$ws.Range("E48").Value = "1.94%"
$ws.Range("E49").Value = "-3.57%"
$ws.Range("E50").Value = "-0.01%"
$ws.Range("E51").Value = "-0.01%"
